# "Fruta / hortaliza, semanal"
#
# A new weekly price observation is inserted as row 625 (dated 2023-08-09,
# serial 45147) in the "Poroto verde" subset sheet. All the existing rows
# from the old row 625 onward shift down by one row to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 625 - shifts rows 625.. down to 626..
$ws.Rows.Item(625).Insert()

# Populate the new row with the new observation.
$ws.Cells.Item(625, 1).Value2 = 3
$ws.Cells.Item(625, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(625, 3).Value2 = "Coquimbo"
$ws.Cells.Item(625, 4).Value2 = 45147
$ws.Cells.Item(625, 5).Value2 = 5
$ws.Cells.Item(625, 6).Value2 = 100112031
$ws.Cells.Item(625, 7).Value2 = "Poroto verde"
$ws.Cells.Item(625, 8).Value2 = "Sin especificar"
$ws.Cells.Item(625, 9).Value2 = "Primera"
$ws.Cells.Item(625, 10).Value2 = 70
$ws.Cells.Item(625, 11).Value2 = 37000
$ws.Cells.Item(625, 12).Value2 = 38000
$ws.Cells.Item(625, 13).Value2 = 37500
$ws.Cells.Item(625, 14).Value2 = "`$/malla 25 kilos"
$ws.Cells.Item(625, 15).Value2 = "Perú"
$ws.Cells.Item(625, 16).Value2 = 1500
$ws.Cells.Item(625, 17).Value2 = 25
$ws.Cells.Item(625, 18).Value2 = "Hortaliza"

# Match the date-serial number format used by the rest of column D.
$ws.Cells.Item(625, 4).NumberFormat = $ws.Cells.Item(626, 4).NumberFormat
